# Scheduled market-data refresh for the Zeromus_Profits leve-crafting sheets.
# Updates the Universalis price snapshot columns (H:N) for a handful of leves
# whose market prices moved since the previous run. Columns:
#   H currentAveragePrice, I currentAveragePriceNQ, J currentAveragePriceHQ,
#   K LevePriceNQ, L LevePriceHQ, M LeveProfitNQ, N LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40: Stuck in the Moment / Horn Glue
$ws.Range("H40").Value = 3356.0356
$ws.Range("I40").Value = 4068.3076
$ws.Range("J40").Value = 2738.7334
$ws.Range("K40").Value = 4068.3076
$ws.Range("L40").Value = 2738.7334
$ws.Range("M40").Value = -3893.3076
$ws.Range("N40").Value = -3088.7334

# Row 46: Always Have an Exit Plan / Poisoning Potion
$ws.Range("H46").Value = 916.6667
$ws.Range("I46").Value = 900
$ws.Range("K46").Value = 2700
$ws.Range("M46").Value = -2581

# Row 60: Make Up Your Mind or Else / Potent Poisoning Potion
$ws.Range("H60").Value = 916.6667
$ws.Range("I60").Value = 900
$ws.Range("K60").Value = 2700
$ws.Range("M60").Value = -2216

# Row 62: The Mustache Suits Him / Enchanted Mythrite Ink
$ws.Range("H62").Value = 3525.75
$ws.Range("J62").Value = 4003
$ws.Range("L62").Value = 4003
$ws.Range("N62").Value = -5251

# Row 64: Forged from the Void / Void Glue
$ws.Range("H64").Value = 2812.724
$ws.Range("I64").Value = 2841.818
$ws.Range("J64").Value = 2794.9443
$ws.Range("K64").Value = 2841.818
$ws.Range("L64").Value = 2794.9443
$ws.Range("M64").Value = -2593.818
$ws.Range("N64").Value = -3290.9443

# Row 65: Forgery of Convenience (L) / Enchanted Mythrite Ink
$ws.Range("H65").Value = 3525.75
$ws.Range("J65").Value = 4003
$ws.Range("L65").Value = 20015
$ws.Range("N65").Value = -26255

# Row 67: Dodging the Draft (L) / Void Glue
$ws.Range("H67").Value = 2812.724
$ws.Range("I67").Value = 2841.818
$ws.Range("J67").Value = 2794.9443
$ws.Range("K67").Value = 2841.818
$ws.Range("L67").Value = 2794.9443
$ws.Range("M67").Value = -1983.818
$ws.Range("N67").Value = -4510.9443

# Row 76: Warding Off Temptation / Enchanted Hardsilver Ink
$ws.Range("H76").Value = 374037.5
$ws.Range("I76").Value = 2222222
$ws.Range("J76").Value = 4400.6
$ws.Range("K76").Value = 2222222
$ws.Range("L76").Value = 4400.6
$ws.Range("M76").Value = -2221907
$ws.Range("N76").Value = -5030.6

# Row 79: The Garden of Arcane Delights (L) / Enchanted Hardsilver Ink
$ws.Range("H79").Value = 374037.5
$ws.Range("I79").Value = 2222222
$ws.Range("J79").Value = 4400.6
$ws.Range("K79").Value = 2222222
$ws.Range("L79").Value = 4400.6
$ws.Range("M79").Value = -2221130
$ws.Range("N79").Value = -6584.6

# Row 98: The Dotted Line / Enchanted Durium Ink
$ws.Range("H98").Value = 774.7778
$ws.Range("I98").Value = 673.75
$ws.Range("J98").Value = 1583
$ws.Range("K98").Value = 673.75
$ws.Range("L98").Value = 1583
$ws.Range("M98").Value = 824.25
$ws.Range("N98").Value = -4579

# Row 112: Making Ends Meet / Superior Spiritbond Potion
$ws.Range("H112").Value = 1028.8846
$ws.Range("J112").Value = 1079.7273
$ws.Range("L112").Value = 3239.1819
$ws.Range("N112").Value = -5455.1819

# Row 122: Wishful Inking / Enchanted High Durium Ink
$ws.Range("H122").Value = 774.7778
$ws.Range("I122").Value = 673.75
$ws.Range("J122").Value = 1583
$ws.Range("K122").Value = 2021.25
$ws.Range("L122").Value = 4749
$ws.Range("M122").Value = 428.75
$ws.Range("N122").Value = -9649

# Row 129: Practical Command / Commanding Craftsman's Draught
$ws.Range("H129").Value = 6372.9487
$ws.Range("I129").Value = 400.76923
$ws.Range("J129").Value = 9359.038
$ws.Range("K129").Value = 1202.30769
$ws.Range("L129").Value = 28077.114
$ws.Range("M129").Value = 3797.69231
$ws.Range("N129").Value = -38077.114

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2450.9136
$ws.Range("J138").Value = 3028.83
$ws.Range("L138").Value = 9086.49
$ws.Range("N138").Value = -19366.49

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 17309.244
$ws.Range("I32").Value = 4833.8237
$ws.Range("J32").Value = 32105.674
$ws.Range("K32").Value = 4833.8237
$ws.Range("L32").Value = 32105.674
$ws.Range("M32").Value = -4546.8237
$ws.Range("N32").Value = -32679.674

# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 2130.92
$ws.Range("I61").Value = 1696.4
$ws.Range("J61").Value = 2782.7
$ws.Range("K61").Value = 1696.4
$ws.Range("L61").Value = 2782.7
$ws.Range("M61").Value = -1484.4
$ws.Range("N61").Value = -3206.7

# Row 110: Scheduled Maintenance / Deepgold Ingot
$ws.Range("H110").Value = 3130.3
$ws.Range("I110").Value = 2301.8333
$ws.Range("J110").Value = 4373
$ws.Range("K110").Value = 2301.8333
$ws.Range("L110").Value = 4373
$ws.Range("M110").Value = -256.8332999999998
$ws.Range("N110").Value = -8463

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 2130.92
$ws.Range("I136").Value = 1696.4
$ws.Range("J136").Value = 2782.7
$ws.Range("K136").Value = 5089.200000000001
$ws.Range("L136").Value = 8348.099999999999
$ws.Range("M136").Value = -2539.200000000001
$ws.Range("N136").Value = -13448.1

$ws = $wb.Worksheets.Item("BSM")
# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 2078.3333
$ws.Range("I107").Value = 1582.5
$ws.Range("J107").Value = 3070
$ws.Range("K107").Value = 1582.5
$ws.Range("L107").Value = 3070
$ws.Range("M107").Value = 337.5
$ws.Range("N107").Value = -6910

$ws = $wb.Worksheets.Item("CRP")
# Row 105: Zelkova, My Love / Zelkova Lumber
$ws.Range("H105").Value = 963.375
$ws.Range("I105").Value = 946
$ws.Range("J105").Value = 992.3333
$ws.Range("K105").Value = 946
$ws.Range("L105").Value = 992.3333
$ws.Range("M105").Value = 801
$ws.Range("N105").Value = -4486.3333

$ws = $wb.Worksheets.Item("CUL")
# Row 107: Slippery Service / Frantoio Oil
$ws.Range("H107").Value = 294.125
$ws.Range("I107").Value = 150
$ws.Range("J107").Value = 342.16666
$ws.Range("K107").Value = 450
$ws.Range("L107").Value = 1026.49998
$ws.Range("M107").Value = 1470
$ws.Range("N107").Value = -4866.499980000001

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 987.1515000000001
$ws.Range("J131").Value = 1044.5
$ws.Range("L131").Value = 3133.5
$ws.Range("N131").Value = -13213.5

$ws = $wb.Worksheets.Item("GSM")
# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 1803.7273
$ws.Range("I122").Value = 2003.6666
$ws.Range("J122").Value = 904
$ws.Range("K122").Value = 6010.9998
$ws.Range("L122").Value = 2712
$ws.Range("M122").Value = -3560.9998
$ws.Range("N122").Value = -7612

# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 2475.4375
$ws.Range("I132").Value = 1744.5
$ws.Range("J132").Value = 3693.6667
$ws.Range("K132").Value = 5233.5
$ws.Range("L132").Value = 11081.0001
$ws.Range("M132").Value = -2703.5
$ws.Range("N132").Value = -16141.0001

$ws = $wb.Worksheets.Item("LTW")
# Row 55: It's Not a Job, It's a Calling / Peiste Leather
$ws.Range("H55").Value = 231.55556
$ws.Range("I55").Value = 232.92308
$ws.Range("J55").Value = 230.28572
$ws.Range("K55").Value = 232.92308
$ws.Range("L55").Value = 230.28572
$ws.Range("M55").Value = -59.92308
$ws.Range("N55").Value = -576.28572

# Row 108: Girding for Glory / Smilodonskin Trousers of Maiming
$ws.Range("H108").Value = 85000
$ws.Range("J108").Value = 85000
$ws.Range("L108").Value = 85000
$ws.Range("N108").Value = -92680

$ws = $wb.Worksheets.Item("WVR")
# Row 75: Storm upon Bald Mountain / Ramie Turban of Crafting
$ws.Range("H75").Value = 24253.334
$ws.Range("J75").Value = 24253.334
$ws.Range("L75").Value = 24253.334
$ws.Range("N75").Value = -26125.334

# Row 78: Abrupt Apprentices (L) / Ramie Turban of Crafting
$ws.Range("H78").Value = 24253.334
$ws.Range("J78").Value = 24253.334
$ws.Range("L78").Value = 72760.00199999999
$ws.Range("N78").Value = -82120.00199999999
